$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "23.899.68"
$ws.Range("E2").Value = "  -1.53%  "
$ws.Range("D3").Value = "1.650.06"
$ws.Range("E3").Value = "  -1.06%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.000"
$ws.Range("E4").Value = "  -0.21%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "310.93"
$ws.Range("E5").Value = "  -0.28%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.3889"
$ws.Range("E7").Value = "  -1.81%  "
$ws.Range("E8").Value = "  -2.50%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "51.18"
$ws.Range("E9").Value = "  -1.39%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "1.344"
$ws.Range("E10").Value = "  -3.00%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "1.000"
$ws.Range("E11").Value = "  -0.20%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.08440"
$ws.Range("E12").Value = "  -1.50%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "23.86"
$ws.Range("E13").Value = "  -2.22%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "7.017"
$ws.Range("E14").Value = "  -3.90%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "7.986"
$ws.Range("E15").Value = "  -0.65%  "
$ws.Range("E16").Value = "  -1.18%  "
$ws.Range("D17").Value = "1.648.24"
$ws.Range("E17").Value = "  -1.25%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "94.03"
$ws.Range("E18").Value = "  -1.79%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.06971"
$ws.Range("E19").Value = "  -0.57%  "
$ws.Range("E20").Value = "  -4.28%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "6.953"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "1.001"
$ws.Range("E22").Value = "  -0.14%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "13.64"
$ws.Range("E23").Value = "  -1.76%  "
$ws.Range("D24").Value = "23.890.78"
$ws.Range("E24").Value = "  -1.59%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.441"
$ws.Range("E25").Value = "  -3.93%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "2.927"
$ws.Range("E26").Value = "  -5.57%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "21.96"
$ws.Range("E27").Value = "  -2.38%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "154.12"
$ws.Range("E28").Value = "  -1.86%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "5.389"
$ws.Range("E29").Value = "  +0.59%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "137.05"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "7.730"
$ws.Range("E31").Value = "  -3.36%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "2.484"
$ws.Range("E32").Value = "  -1.80%  "
$ws.Range("D33").Value = "1.833.19"
$ws.Range("E33").Value = "  -0.92%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.08163"
$ws.Range("E34").Value = "  -0.97%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.9947"
$ws.Range("E35").Value = "  -6.05%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.2679"
$ws.Range("E38").Value = "  -3.13%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "10.55"
$ws.Range("E39").Value = "  -5.35%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.09108"
$ws.Range("E40").Value = "  -1.85%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.7562"
$ws.Range("E41").Value = "  -1.47%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "13.51"
$ws.Range("E42").Value = "  -2.00%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "1.422"
$ws.Range("E43").Value = "  -1.61%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "16.74"
$ws.Range("E44").Value = "  +0.56%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.6936"
$ws.Range("E45").Value = "  -1.74%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "2.448"
$ws.Range("E46").Value = "  -3.16%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "4.099"
$ws.Range("E47").Value = "  -0.58%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.9999"
$ws.Range("E48").Value = "  -0.12%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.08266"
$ws.Range("E49").Value = "  -1.61%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "134.08"
$ws.Range("E50").Value = "  -1.76%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "1.224"
$ws.Range("E51").Value = "  -2.75%  "

# Row 36 and 37: coin identities swap (name/link) with refreshed price/volume data
$ws.Range("B36").Value = "InternetComputer(DFINITY)"
$ws.Range("C36").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "6.708"
$ws.Range("E36").Value = "  -2.22%  "
$ws.Range("B37").Value = "VeChain"
$ws.Range("C37").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.02913"
$ws.Range("E37").Value = "  -4.89%  "
